# "Created graphs for Percentage Feature Creators"
# Update the data on Sheet1 (row 3 values) and move the active selection,
# matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 data edits
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 9
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 6

# Move the active cell / selection from B10 to H8
$ws.Range("H8").Select()

# Reposition the workbook window (best-effort; mirrors the author's
# window-position change in the saved view state).
$win = $wb.Windows.Item(1)
$win.Left = 1160
$win.Top = 1460
